# Restored from revision #ba2fc01d22c321458b86684e19066dd142eb40f1.TEST
# The only substantive change in that revision is the "min" threshold for
# the R30 rule (row 10) on the Rules sheet: it goes back from 18 to 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 20
